$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.898.93'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').Value = '2.597.75'
$ws.Range('E3').Value = '  -1.72%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '552.63'
$ws.Range('E5').Value = '  +2.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.97'
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.600'
$ws.Range('E8').Value = '  +4.93%  '
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('E10').Value = '  -1.44%  '
$ws.Range('E11').Value = '  +5.08%  '
$ws.Range('E12').Value = '  -0.36%  '
$ws.Range('D13').Value = '3.053.50'
$ws.Range('E13').Value = '  -1.75%  '
$ws.Range('D14').Value = '58.857.83'
$ws.Range('E14').Value = '  -1.00%  '
$ws.Range('E15').Value = '  -1.41%  '
$ws.Range('D16').Value = '2.595.94'
$ws.Range('E16').Value = '  -2.66%  '
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.07'
$ws.Range('E20').Value = '  -2.77%  '
$ws.Range('E21').Value = '  -1.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.47'
$ws.Range('E23').Value = '  -0.75%  '
$ws.Range('E24').Value = '  +2.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.996'
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.160'
$ws.Range('E26').Value = '  -2.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.13'
$ws.Range('E27').Value = '  -1.89%  '
$ws.Range('E28').Value = '  +1.70%  '
$ws.Range('E30').Value = '  +1.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '152.91'
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.96'
$ws.Range('E33').Value = '  +0.29%  '
$ws.Range('E34').Value = '  -2.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.877'
$ws.Range('E35').Value = '  +3.63%  '
$ws.Range('E36').Value = '  -1.33%  '
$ws.Range('E37').Value = '  -0.52%  '
$ws.Range('E38').Value = '  +1.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.824'
$ws.Range('E39').Value = '  -1.26%  '
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '284.01'
$ws.Range('E41').Value = '  -1.37%  '
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.598'
$ws.Range('E43').Value = '  -1.16%  '
$ws.Range('E44').Value = '  +1.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.61'
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('E46').Value = '  -0.94%  '
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '119.48'
$ws.Range('E48').Value = '  +7.84%  '
$ws.Range('D49').Value = '1.923.41'
$ws.Range('E49').Value = '  -2.35%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.43'
$ws.Range('E50').Value = '  -2.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.85'
$ws.Range('E51').Value = '  -2.45%  '
